$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 3431.5625
$ws.Cells.Item(40, 9).Value = 3300.2
$ws.Cells.Item(40, 10).Value = 3491.2727
$ws.Cells.Item(40, 11).Value = 3300.2
$ws.Cells.Item(40, 12).Value = 3491.2727
$ws.Cells.Item(40, 13).Value = -3125.2
$ws.Cells.Item(40, 14).Value = -3841.2727
$ws.Cells.Item(76, 8).Value = 3625
$ws.Cells.Item(76, 9).Value = 3400
$ws.Cells.Item(76, 10).Value = 3850
$ws.Cells.Item(76, 11).Value = 3400
$ws.Cells.Item(76, 12).Value = 3850
$ws.Cells.Item(76, 13).Value = -3085
$ws.Cells.Item(76, 14).Value = -4480
$ws.Cells.Item(79, 8).Value = 3625
$ws.Cells.Item(79, 9).Value = 3400
$ws.Cells.Item(79, 10).Value = 3850
$ws.Cells.Item(79, 11).Value = 3400
$ws.Cells.Item(79, 12).Value = 3850
$ws.Cells.Item(79, 13).Value = -2308
$ws.Cells.Item(79, 14).Value = -6034
$ws.Cells.Item(132, 8).Value = 7449.6978
$ws.Cells.Item(132, 9).Value = 8511.190000000001
$ws.Cells.Item(132, 10).Value = 6436.4546
$ws.Cells.Item(132, 11).Value = 25533.57
$ws.Cells.Item(132, 12).Value = 19309.3638
$ws.Cells.Item(132, 13).Value = -23003.57
$ws.Cells.Item(132, 14).Value = -24369.3638
$ws.Cells.Item(135, 8).Value = 543.9259
$ws.Cells.Item(135, 9).Value = 336.91666
$ws.Cells.Item(135, 10).Value = 2200
$ws.Cells.Item(135, 11).Value = 3032.24994
$ws.Cells.Item(135, 12).Value = 19800
$ws.Cells.Item(135, 13).Value = -497.2499399999997
$ws.Cells.Item(135, 14).Value = -24870
$ws.Cells.Item(138, 8).Value = 1643.8853
$ws.Cells.Item(138, 9).Value = 1118.1464
$ws.Cells.Item(138, 10).Value = 2721.65
$ws.Cells.Item(138, 11).Value = 3354.4392
$ws.Cells.Item(138, 12).Value = 8164.950000000001
$ws.Cells.Item(138, 13).Value = 1785.5608
$ws.Cells.Item(138, 14).Value = -18444.95
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 7191.6885
$ws.Cells.Item(32, 9).Value = 7194.891
$ws.Cells.Item(32, 11).Value = 7194.891
$ws.Cells.Item(32, 13).Value = -6907.891
$ws.Cells.Item(61, 8).Value = 2369.318
$ws.Cells.Item(61, 9).Value = 2910.375
$ws.Cells.Item(61, 10).Value = 2060.1428
$ws.Cells.Item(61, 11).Value = 2910.375
$ws.Cells.Item(61, 12).Value = 2060.1428
$ws.Cells.Item(61, 13).Value = -2698.375
$ws.Cells.Item(61, 14).Value = -2484.1428
$ws.Cells.Item(136, 8).Value = 2369.318
$ws.Cells.Item(136, 9).Value = 2910.375
$ws.Cells.Item(136, 10).Value = 2060.1428
$ws.Cells.Item(136, 11).Value = 8731.125
$ws.Cells.Item(136, 12).Value = 6180.428400000001
$ws.Cells.Item(136, 13).Value = -6181.125
$ws.Cells.Item(136, 14).Value = -11280.4284
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 2179.1482
$ws.Cells.Item(99, 9).Value = 1628.85
$ws.Cells.Item(99, 10).Value = 3751.4285
$ws.Cells.Item(99, 11).Value = 1628.85
$ws.Cells.Item(99, 12).Value = 3751.4285
$ws.Cells.Item(99, 13).Value = -130.8499999999999
$ws.Cells.Item(99, 14).Value = -6747.4285
$ws.Cells.Item(134, 8).Value = 5449.1665
$ws.Cells.Item(134, 9).Value = 2779.0557
$ws.Cells.Item(134, 10).Value = 7451.75
$ws.Cells.Item(134, 11).Value = 8337.167099999999
$ws.Cells.Item(134, 12).Value = 22355.25
$ws.Cells.Item(134, 13).Value = -5802.167099999999
$ws.Cells.Item(134, 14).Value = -27425.25
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 9011583
$ws.Cells.Item(31, 9).Value = 1563.625
$ws.Cells.Item(31, 10).Value = 25645464
$ws.Cells.Item(31, 11).Value = 1563.625
$ws.Cells.Item(31, 12).Value = 25645464
$ws.Cells.Item(31, 13).Value = -1268.625
$ws.Cells.Item(31, 14).Value = -25646054
$ws.Cells.Item(34, 8).Value = 9011583
$ws.Cells.Item(34, 9).Value = 1563.625
$ws.Cells.Item(34, 10).Value = 25645464
$ws.Cells.Item(34, 11).Value = 1563.625
$ws.Cells.Item(34, 12).Value = 25645464
$ws.Cells.Item(34, 13).Value = -1361.625
$ws.Cells.Item(34, 14).Value = -25645868
$ws.Cells.Item(132, 8).Value = 3595.75
$ws.Cells.Item(132, 9).Value = 3288
$ws.Cells.Item(132, 10).Value = 3749.625
$ws.Cells.Item(132, 11).Value = 9864
$ws.Cells.Item(132, 12).Value = 11248.875
$ws.Cells.Item(132, 13).Value = -7334
$ws.Cells.Item(132, 14).Value = -16308.875
$ws.Cells.Item(134, 8).Value = 1997.5
$ws.Cells.Item(134, 9).Value = 1339.3334
$ws.Cells.Item(134, 10).Value = 2177
$ws.Cells.Item(134, 11).Value = 4018.0002
$ws.Cells.Item(134, 12).Value = 6531
$ws.Cells.Item(134, 13).Value = -1483.0002
$ws.Cells.Item(134, 14).Value = -11601
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(88, 8).Value = 4356.3125
$ws.Cells.Item(88, 10).Value = 4356.3125
$ws.Cells.Item(88, 12).Value = 13068.9375
$ws.Cells.Item(88, 14).Value = -13924.9375
$ws.Cells.Item(91, 8).Value = 4356.3125
$ws.Cells.Item(91, 10).Value = 4356.3125
$ws.Cells.Item(91, 12).Value = 13068.9375
$ws.Cells.Item(91, 14).Value = -16032.9375
$ws.Cells.Item(93, 8).Value = 5313.478
$ws.Cells.Item(93, 10).Value = 5313.478
$ws.Cells.Item(93, 12).Value = 15940.434
$ws.Cells.Item(93, 14).Value = -19684.434
$ws.Cells.Item(98, 8).Value = 1151
$ws.Cells.Item(98, 10).Value = 2322.1667
$ws.Cells.Item(98, 12).Value = 6966.500100000001
$ws.Cells.Item(98, 14).Value = -9962.500100000001
$ws.Cells.Item(105, 8).Value = 5947.4546
$ws.Cells.Item(105, 10).Value = 6039.6
$ws.Cells.Item(105, 12).Value = 18118.8
$ws.Cells.Item(105, 14).Value = -23360.8
$ws.Cells.Item(121, 8).Value = 2312
$ws.Cells.Item(121, 9).Value = 277.14285
$ws.Cells.Item(121, 10).Value = 3407.6924
$ws.Cells.Item(121, 11).Value = 831.4285500000001
$ws.Cells.Item(121, 12).Value = 10223.0772
$ws.Cells.Item(121, 13).Value = 478.5714499999999
$ws.Cells.Item(121, 14).Value = -12843.0772
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5813.5864
$ws.Cells.Item(70, 9).Value = 4480.6
$ws.Cells.Item(70, 10).Value = 6515.1577
$ws.Cells.Item(70, 11).Value = 4480.6
$ws.Cells.Item(70, 12).Value = 6515.1577
$ws.Cells.Item(70, 13).Value = -4210.6
$ws.Cells.Item(70, 14).Value = -7055.1577
$ws.Cells.Item(73, 8).Value = 5813.5864
$ws.Cells.Item(73, 9).Value = 4480.6
$ws.Cells.Item(73, 10).Value = 6515.1577
$ws.Cells.Item(73, 11).Value = 4480.6
$ws.Cells.Item(73, 12).Value = 6515.1577
$ws.Cells.Item(73, 13).Value = -3544.6
$ws.Cells.Item(73, 14).Value = -8387.1577
$ws.Cells.Item(80, 8).Value = 3031.818
$ws.Cells.Item(80, 9).Value = 2964.375
$ws.Cells.Item(80, 10).Value = 3070.3572
$ws.Cells.Item(80, 11).Value = 2964.375
$ws.Cells.Item(80, 12).Value = 3070.3572
$ws.Cells.Item(80, 13).Value = -1966.375
$ws.Cells.Item(80, 14).Value = -5066.3572
$ws.Cells.Item(83, 8).Value = 3031.818
$ws.Cells.Item(83, 9).Value = 2964.375
$ws.Cells.Item(83, 10).Value = 3070.3572
$ws.Cells.Item(83, 11).Value = 14821.875
$ws.Cells.Item(83, 12).Value = 15351.786
$ws.Cells.Item(83, 13).Value = -9829.875
$ws.Cells.Item(83, 14).Value = -25335.786
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 3579.9
$ws.Cells.Item(100, 9).Value = 3300
$ws.Cells.Item(100, 10).Value = 3859.8
$ws.Cells.Item(100, 11).Value = 3300
$ws.Cells.Item(100, 12).Value = 3859.8
$ws.Cells.Item(100, 13).Value = -2759
$ws.Cells.Item(100, 14).Value = -4941.8
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1832.9592
$ws.Cells.Item(132, 9).Value = 1327.6666
$ws.Cells.Item(132, 10).Value = 2875.125
$ws.Cells.Item(132, 11).Value = 3982.9998
$ws.Cells.Item(132, 12).Value = 8625.375
$ws.Cells.Item(132, 13).Value = -1452.9998
$ws.Cells.Item(132, 14).Value = -13685.375

Write-Host "Applied 171 cell updates across 8 sheets"